# Updated cryptos list with GitHub Actions scrape refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.962.91"
$ws.Range("E2").Value = "  +2.09%  "
$ws.Range("D3").Value = "'3.739.25"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'601.80"
$ws.Range("E5").Value = "  +1.59%  "
$ws.Range("D6").Value = "'168.06"
$ws.Range("E6").Value = "  -1.90%  "
$ws.Range("D7").Value = "'3.739.08"
$ws.Range("E7").Value = "  +0.32%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +3.33%  "
$ws.Range("E10").Value = "  +5.35%  "
$ws.Range("D11").Value = "'6.34"
$ws.Range("E11").Value = "  +2.79%  "
$ws.Range("E12").Value = "  +0.39%  "
$ws.Range("D13").Value = "'38.18"
$ws.Range("E13").Value = "  +2.01%  "
$ws.Range("E14").Value = "  +1.75%  "
$ws.Range("D15").Value = "'4.362.86"
$ws.Range("E15").Value = "  +0.30%  "
$ws.Range("D16").Value = "'3.735.69"
$ws.Range("E16").Value = "  +0.39%  "
$ws.Range("D17").Value = "'68.911.47"
$ws.Range("E17").Value = "  +2.01%  "
$ws.Range("D18").Value = "'7.25"
$ws.Range("E18").Value = "  +1.82%  "
$ws.Range("E19").Value = "  +0.99%  "
$ws.Range("D20").Value = "'17.23"
$ws.Range("E20").Value = "  +6.84%  "
$ws.Range("D21").Value = "'496.85"
$ws.Range("E21").Value = "  +2.00%  "
$ws.Range("D22").Value = "'10.09"
$ws.Range("E22").Value = "  +12.46%  "
$ws.Range("E23").Value = "  +1.46%  "
$ws.Range("D24").Value = "'85.28"
$ws.Range("E24").Value = "  +2.36%  "
$ws.Range("E25").Value = "  -1.09%  "
$ws.Range("E26").Value = "  -0.54%  "
$ws.Range("D27").Value = "'12.31"
$ws.Range("E27").Value = "  +1.74%  "
$ws.Range("E28").Value = "  -0.51%  "
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("E30").Value = "  +1.00%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").Value = "'2.42"
$ws.Range("E31").Value = "  +2.21%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").Value = "'7.96"
$ws.Range("E32").Value = "  +4.50%  "
$ws.Range("D33").Value = "'31.74"
$ws.Range("E33").Value = "  -1.39%  "
$ws.Range("D34").Value = "'3.889.70"
$ws.Range("E34").Value = "  +0.63%  "
$ws.Range("E35").Value = "  +0.82%  "
$ws.Range("D36").Value = "'3.669.57"
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("E38").Value = "  +1.57%  "
$ws.Range("E39").Value = "  +2.16%  "
$ws.Range("E40").Value = "  -0.25%  "
$ws.Range("E41").Value = "  +0.68%  "
$ws.Range("D42").Value = "'436.23"
$ws.Range("E42").Value = "  -2.42%  "
$ws.Range("D43").Value = "'48.99"
$ws.Range("E43").Value = "  +0.41%  "
$ws.Range("E44").Value = "  +0.74%  "
$ws.Range("D45").Value = "'2.86"
$ws.Range("E45").Value = "  +0.79%  "
$ws.Range("E46").Value = "  +2.13%  "
$ws.Range("D48").Value = "'40.50"
$ws.Range("E48").Value = "  -1.71%  "
$ws.Range("D49").Value = "'142.06"
$ws.Range("E49").Value = "  +0.92%  "
$ws.Range("D50").Value = "'0.0353"
$ws.Range("E50").Value = "  +2.11%  "
$ws.Range("D51").Value = "'2.747.67"
$ws.Range("E51").Value = "  -1.25%  "
